$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3, shifting the header row (and everything
# below it) down by one. This turns the former "gap" (row 2, from the
# A1:F2 merge) into a two-row gap (rows 2-3) before the table starts at row 4.
$ws.Rows("3").Insert()
